# Applies the "Updated capital structure database" refresh to the Singapore
# Investments & Asset Management sheet: rows 2-4 get refreshed metrics, a new
# row for "Uni-Asia Group Limited" is inserted after row 4 (shifting TIH Limited,
# Reenova->TIH, 8i Enterprises->SC Health down), Reenova Investment Holding is
# dropped, and the final row becomes "Forise International Limited".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 0.127
$ws.Range("E2").Value = 0.112
$ws.Range("F2").Value = 0.388
$ws.Range("G2").Value = 0.2023340765690972
$ws.Range("H2").Value = 0.1895251197874662
$ws.Range("I2").Value = 0.07391637806980091
$ws.Range("J2").Value = 0.07174236695010087
$ws.Range("K2").Value = -1.326000000000001
$ws.Range("L2").Value = -0.005242184164334175
$ws.Range("M2").Value = 10.457
$ws.Range("N2").Value = 0.01043998282799038
$ws.Range("O2").Value = -7.886123680241323
$ws.Range("P2").Value = 10.09
$ws.Range("Q2").Value = 0.01007358006449487
$ws.Range("R2").Value = -7.609351432880841
$ws.Range("S2").Value = 0.367
$ws.Range("T2").Value = 0.0350961078703261
$ws.Range("U2").Value = 69.37400000000001
$ws.Range("V2").Value = 0.06926110439982829
$ws.Range("W2").Value = -0.04072845838800747
$ws.Range("X2").Value = 0.03536307638425101
$ws.Range("Y2").Value = -0.07609153477225848
$ws.Range("Z2").Value = 0.3768189347984125
$ws.Range("AA2").Value = -0.01807167951781926
$ws.Range("AB2").Value = 0.035279085572385
$ws.Range("AC2").Value = -0.05420242512310176
$ws.Range("AD2").Value = 153.535
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 153.535
$ws.Range("AG2").Value = 84.16099999999999
$ws.Range("AH2").Value = 0.1329117485380876
$ws.Range("AI2").Value = 0.3320429502914175
$ws.Range("AJ2").Value = 0.07751123374572086
$ws.Range("AK2").Value = 0.2141386847013264
$ws.Range("AL2").Value = 10.152
$ws.Range("AM2").Value = 7.570000000000001
$ws.Range("AN2").Value = 4.664023816033293
$ws.Range("AO2").Value = 1.841706067769897
$ws.Range("AP2").Value = 2.556608645463106
$ws.Range("AQ2").Value = 2.469881109643329

# --- Row 3 ---
$ws.Range("G3").Value = 0.3633663366336634
$ws.Range("H3").Value = 0.3633663366336634
$ws.Range("I3").Value = 0.2712871287128713
$ws.Range("J3").Value = 0.2712871287128713
$ws.Range("K3").Value = -1.8
$ws.Range("L3").Value = -0.1782178217821782
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("T3").Value = $null
$ws.Range("X3").Value = 0.03530966097083269
$ws.Range("Z3").Value = $null
$ws.Range("AA3").Value = $null
$ws.Range("AB3").Value = 0.03525210685652549
$ws.Range("AC3").Value = $null
$ws.Range("AD3").Value = 0.219
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.219
$ws.Range("AG3").Value = 0.219
$ws.Range("AH3").Value = 0.00297074024335653
$ws.Range("AI3").Value = 0.01017705283702774
$ws.Range("AJ3").Value = 0.00297074024335653
$ws.Range("AK3").Value = 0.01017705283702774
$ws.Range("AL3").Value = 0.014
$ws.Range("AM3").Value = 0.014
$ws.Range("AN3").Value = 0.07577854671280276
$ws.Range("AO3").Value = 195.7142857142857
$ws.Range("AP3").Value = 0.07577854671280276
$ws.Range("AQ3").Value = 195.7142857142857

# --- Row 4 ---
$ws.Range("D4").Value = 0.131
$ws.Range("E4").Value = 0.112
$ws.Range("F4").Value = 0.388
$ws.Range("G4").Value = 0.1014886164623468
$ws.Range("H4").Value = 0.07311733800350262
$ws.Range("I4").Value = 0.1208406304728546
$ws.Range("J4").Value = 0.09951581333058618
$ws.Range("K4").Value = 12.7
$ws.Range("L4").Value = 0.1112084063047285
$ws.Range("M4").Value = 6.598
$ws.Range("N4").Value = 0.01067119521267993
$ws.Range("O4").Value = 0.5195275590551182
$ws.Range("P4").Value = 6.25
$ws.Range("Q4").Value = 0.01010836163674592
$ws.Range("R4").Value = 0.4921259842519685
$ws.Range("S4").Value = 0.3479999999999999
$ws.Range("T4").Value = 0.05274325553197937
$ws.Range("U4").Value = 23.8
$ws.Range("V4").Value = 0.03849264111272845
$ws.Range("W4").Value = 0.1978193146417445
$ws.Range("X4").Value = 0.03584803762187968
$ws.Range("Y4").Value = 0.1619712770198648
$ws.Range("Z4").Value = 1.691851851851851
$ws.Range("AA4").Value = 0.1683660130718954
$ws.Range("AB4").Value = 0.0353060642882445
$ws.Range("AC4").Value = 0.1330599487836509
$ws.Range("AD4").Value = 17.3
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 17.3
$ws.Range("AG4").Value = -6.5
$ws.Range("AH4").Value = 0.02721837633731907
$ws.Range("AI4").Value = 0.1911602209944751
$ws.Range("AJ4").Value = -0.010624387054593
$ws.Range("AK4").Value = -0.09745127436281859
$ws.Range("AL4").Value = 0.643
$ws.Range("AM4").Value = -0.5469999999999999
$ws.Range("AN4").Value = 1.074534161490683
$ws.Range("AO4").Value = 21.46189735614308
$ws.Range("AP4").Value = -0.4037267080745341
$ws.Range("AQ4").Value = -25.22851919561244

# --- Row 5 ---
$ws.Range("B5").Value = "Uni-Asia Group Limited (SGX:CHJ)"
$ws.Range("D5").Value = 0.123
$ws.Range("G5").Value = 0.202020202020202
$ws.Range("H5").Value = 0.202020202020202
$ws.Range("I5").Value = 0.06348096348096349
$ws.Range("J5").Value = 0.06348096348096349
$ws.Range("K5").Value = -4.43
$ws.Range("L5").Value = -0.03442113442113443
$ws.Range("M5").Value = 2.089
$ws.Range("N5").Value = 0.05901129943502825
$ws.Range("O5").Value = -0.4715575620767495
$ws.Range("P5").Value = 2.07
$ws.Range("Q5").Value = 0.05847457627118644
$ws.Range("R5").Value = -0.4672686230248307
$ws.Range("S5").Value = 0.01900000000000013
$ws.Range("T5").Value = 0.009095260890378232
$ws.Range("U5").Value = 31.1
$ws.Range("V5").Value = 0.8785310734463277
$ws.Range("W5").Value = -0.03558232931726907
$ws.Range("X5").Value = 0.1179170818065115
$ws.Range("Y5").Value = -0.1534994111237806
$ws.Range("Z5").Value = 0.2503404007002528
$ws.Range("AA5").Value = 0.01589184983466252
$ws.Range("AB5").Value = 0.03701089956772041
$ws.Range("AC5").Value = -0.0211190497330579
$ws.Range("AD5").Value = 135.9
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 135.9
$ws.Range("AG5").Value = 104.8
$ws.Range("AH5").Value = 0.7933450087565674
$ws.Range("AI5").Value = 0.52899961074348
$ws.Range("AJ5").Value = 0.7475035663338089
$ws.Range("AK5").Value = 0.4641275465013286
$ws.Range("AL5").Value = 9.470000000000001
$ws.Range("AM5").Value = 9.470000000000001
$ws.Range("AN5").Value = 7.152631578947369
$ws.Range("AO5").Value = 0.8627243928194297
$ws.Range("AP5").Value = 5.515789473684211
$ws.Range("AQ5").Value = 0.8627243928194297

# --- Row 6 ---
$ws.Range("B6").Value = "TIH Limited (SGX:T55)"
$ws.Range("G6").Value = -20.19047619047619
$ws.Range("H6").Value = -20.19047619047619
$ws.Range("I6").Value = 8.152380952380952
$ws.Range("J6").Value = 8.152380952380952
$ws.Range("K6").Value = -4.17
$ws.Range("L6").Value = 7.942857142857142
$ws.Range("M6").Value = 1.77
$ws.Range("N6").Value = 0.04836065573770491
$ws.Range("O6").Value = -0.4244604316546763
$ws.Range("P6").Value = 1.77
$ws.Range("Q6").Value = 0.04836065573770491
$ws.Range("R6").Value = -0.4244604316546763
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 13.8
$ws.Range("V6").Value = 0.3770491803278689
$ws.Range("W6").Value = -0.04587458745874587
$ws.Range("X6").Value = 0.03529492021382818
$ws.Range("Y6").Value = -0.08116950767257405
$ws.Range("Z6").Value = -0.006382823517969169
$ws.Range("AA6").Value = -0.05203520887030103
$ws.Range("AB6").Value = 0.0352505916428446
$ws.Range("AC6").Value = -0.08728580051314563
$ws.Range("AD6").Value = 0.08400000000000001
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0.08400000000000001
$ws.Range("AG6").Value = -13.716
$ws.Range("AH6").Value = 0.002289826627412496
$ws.Range("AI6").Value = 0.0009735292754160679
$ws.Range("AJ6").Value = -0.5993707393812271
$ws.Range("AK6").Value = -0.1892279675514596
$ws.Range("AM6").Value = -0.052
$ws.Range("AN6").Value = -0.01981132075471698
$ws.Range("AP6").Value = 3.234905660377358
$ws.Range("AQ6").Value = 82.30769230769232

# --- Row 7 ---
$ws.Range("B7").Value = "SC Health Corporation (NYSE:SCPE)"
$ws.Range("K7").Value = 0.444
$ws.Range("O7").Value = -0
$ws.Range("R7").Value = -0
$ws.Range("U7").Value = 0.111
$ws.Range("V7").Value = 0.000474764756201882
$ws.Range("W7").Value = 0.0888
$ws.Range("X7").Value = 0.03524549617030377
$ws.Range("Y7").Value = 0.05355450382969624
$ws.Range("AA7").Value = -0.2278481012658228
$ws.Range("AB7").Value = 0.03524549617030377
$ws.Range("AC7").Value = -0.2630935974361265
$ws.Range("AD7").Value = 0
$ws.Range("AF7").Value = 0
$ws.Range("AG7").Value = -0.111
$ws.Range("AH7").Value = 0
$ws.Range("AI7").Value = 0
$ws.Range("AJ7").Value = -0.0004749902648391666
$ws.Range("AK7").Value = -0.02270402945387605
$ws.Range("AM7").Value = -1.34
$ws.Range("AQ7").Value = 0.6716417910447761

# --- Row 8 ---
$ws.Range("B8").Value = "Forise International Limited (SGX:8A1)"
$ws.Range("D8").Value = $null
$ws.Range("E8").Value = $null
$ws.Range("G8").Value = -1.437632135306554
$ws.Range("H8").Value = -1.437632135306554
$ws.Range("I8").Value = -1.761099365750529
$ws.Range("J8").Value = -1.761099365750529
$ws.Range("K8").Value = -4.07
$ws.Range("L8").Value = -8.604651162790699
$ws.Range("M8").Value = -0
$ws.Range("N8").Value = -0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = -0
$ws.Range("Q8").Value = -0
$ws.Range("R8").Value = 0
$ws.Range("T8").Value = $null
$ws.Range("U8").Value = 0.5629999999999999
$ws.Range("V8").Value = 0.1397022332506203
$ws.Range("W8").Value = -0.6339563862928349
$ws.Range("X8").Value = 0.03541649179766934
$ws.Range("Y8").Value = -0.6693728780905043
$ws.Range("Z8").Value = 0.1363112391930836
$ws.Range("AA8").Value = -0.2400576368876081
$ws.Range("AB8").Value = 0.03533351323107027
$ws.Range("AC8").Value = -0.2753911501186784
$ws.Range("AD8").Value = 0.032
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 0.032
$ws.Range("AG8").Value = -0.5309999999999999
$ws.Range("AH8").Value = 0.007877892663712456
$ws.Range("AI8").Value = 0.0145985401459854
$ws.Range("AJ8").Value = -0.1517576450414404
$ws.Range("AK8").Value = -0.3259668508287292
$ws.Range("AL8").Value = 0.025
$ws.Range("AM8").Value = 0.025
$ws.Range("AN8").Value = -0.03850782190132371
$ws.Range("AO8").Value = -33.31999999999999
$ws.Range("AP8").Value = 0.6389891696750902
$ws.Range("AQ8").Value = -33.31999999999999
